# Update column F (dSF) values for rows 2, 4, 5, 17, 18 per repulled data / mean calculation
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 6
$ws.Range("F4").Value = 7
$ws.Range("F5").Value = -3
$ws.Range("F17").Value = -7
$ws.Range("F18").Value = -2
